# Apply "added download history of leases" edit.
#
# The sheet previously held one row per currently-active lease (with an
# "active" boolean flag, plus slug / shouldBeReturned / id bookkeeping
# columns). It is replaced by a lease *download/history* export: the
# lease id moves into column B ("orderedBook"), the "active" flag and the
# slug/shouldBeReturned/id columns are dropped, and the row set/order is
# refreshed to the new export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "studentName"
$ws.Range("B1").Value = "orderedBook"
$ws.Range("C1").Value = "orderedBookSeria"
$ws.Range("D1").Value = "classOfStudent"
$ws.Range("E1").Value = "major"
$ws.Range("F1").Value = "studentPhoneNumber"
$ws.Range("G1").Value = "orderedTime"
$ws.Range("H1").Value = "deadline"

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = "Bakhodirova  Rukhsorakhon  Behzod qizi"
$ws.Range("B2").Value = "6389c809666ce69a7f92ffa9"
$ws.Range("C2").Value = "CH-00118"
$ws.Range("D2").Value = "E-203"
$ws.Range("E2").Value = "ECE"
$ws.Range("F2").Value = " 998 97 4633882"
$ws.Range("G2").Value = "2022-12-10T09:36:53.998Z"
$ws.Range("H2").Value = "2022-12-17T09:36:53.998Z"

# --- Row 3 --------------------------------------------------------------
# Name contains U+2018 (LEFT SINGLE QUOTATION MARK) twice: o{u2018}g{u2018}li
$ws.Range("A3").Value = "Tulaev  Muhammad  Dilshod o" + [char]0x2018 + "g" + [char]0x2018 + "li"
$ws.Range("B3").Value = "638f8d2f008bc82d3ca29407"
$ws.Range("C3").Value = "N-01597"
$ws.Range("D3").Value = "E-203"
$ws.Range("E3").Value = "ECE"
$ws.Range("F3").Value = " 998 90 864-20-24"
$ws.Range("G3").Value = "2022-12-10T09:33:35.871Z"
$ws.Range("H3").Value = "2022-12-17T09:33:35.871Z"

# --- Row 4 ---------------------------------------------------------------
# A4/C4/D4/E4/F4 are untouched by the diff (already correct); only G4/H4
# change shape (boolean "active" flag -> plain orderedTime string, and the
# deadline shifts left from column I into H).
$ws.Range("G4").Value = "2022-12-06T21:03:48.510Z"
$ws.Range("H4").Value = "2022-12-13T21:03:48.510Z"

# --- Row 5 --------------------------------------------------------------
$ws.Range("A5").Value = "Mahmudhodjayev Izzat"
$ws.Range("B5").Value = "6389c809666ce69a7f92ffaa"
$ws.Range("C5").Value = "CH-00017"
$ws.Range("D5").Value = "A203"
$ws.Range("E5").Value = "Architecture"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "+998977654321"
$ws.Range("G5").Value = "2022-12-06T11:35:38.491Z"
$ws.Range("H5").Value = "2022-12-13T11:35:38.491Z"

# --- Row 6 --------------------------------------------------------------
$ws.Range("A6").Value = "Jo'rayev Narimon"
$ws.Range("B6").Value = "6389c809666ce69a7f92ffac"
$ws.Range("C6").Value = "CH-00092"
$ws.Range("D6").Value = "E202"
$ws.Range("E6").Value = "ECE"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "+998911234567"
$ws.Range("G6").Value = "2022-12-06T11:19:56.130Z"
$ws.Range("H6").Value = "2022-12-13T11:19:56.130Z"

# --- Drop the old tail: extra columns I:L and the old 7th row -----------
$ws.Range("I1:L7").Clear() | Out-Null
$ws.Range("A7:H7").Clear() | Out-Null
